# Updates "cryptos" worksheet with refreshed price / volume(1h) figures
# (and, for rows 49-50, an Elrond/Cronos content swap) per the latest
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: several "Price" values look numeric (e.g. "13.90", "1.001") but
# must stay as literal text (matching the source feed's formatting,
# incl. trailing zeros). Assigning such strings straight to .Value lets
# Excel "smart type" them into real numbers and silently drop trailing
# zeros/precision, so for those cells we first force the cell to Text
# format ("@") before writing the string.

$ws.Range("D2").Value = "26.596.35"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.842.80"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.58"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3162"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06793"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.13"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7834"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "1.838.73"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.26"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.014"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.90"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007928"
$ws.Range("D20").Value = "26.629.77"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "2.073.52"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.611"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.993"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.332"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.228"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.01"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.87"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.215"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08698"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.079"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04859"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7316"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.864"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.105"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.344"
$ws.Range("E38").Value = "  +5.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01732"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9045"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.26"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.916"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.705"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.064"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1244"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.92"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05824"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8941"
$ws.Range("E51").Value = "  +0.86%  "
